{"js": "// Locate the \"Chapter 8\" paragraph (contains \"#8.16\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet chapter8Para = null;\nfor (const p of paragraphs.items) {\n  if (p.text.indexOf(\"#8.16\") !== -1) {\n    chapter8Para = p;\n    break;\n  }\n}\n\nif (!chapter8Para) {\n  throw new Error(\"Could not find the 'Chapter 8' paragraph containing '#8.16'.\");\n}\n\n// 1) Remove the existing \"_GoBack\" bookmark (it currently sits after \"9.13\"\n//    in the \"Chapter 9\" paragraph) so it can be re-inserted at its new spot.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Split the \"#8.16\" run into \"#8.\" + \"16\" by inserting the \"_GoBack\"\n//    bookmark right before the \"16\".\nconst fullMatches = chapter8Para.search(\"#8.16\", { matchCase: true });\nfullMatches.load(\"items\");\nawait context.sync();\n\nconst fullRange = fullMatches.items[0];\nconst sixteenMatches = fullRange.search(\"16\", { matchCase: true });\nsixteenMatches.load(\"items\");\nawait context.sync();\n\nconst sixteenRange = sixteenMatches.items[0];\nconst insertionPoint = sixteenRange.getRange(\"Start\");\ninsertionPoint.insertBookmark(\"_GoBack\");\nawait context.sync();\n\n// 3) Highlight the \", # 8.20\" text (the 4 runs: \",\", \" #\", \" \", \"8.20\")\n//    in yellow.\nconst highlightMatches = chapter8Para.search(\", # 8.20\", { matchCase: true });\nhighlightMatches.load(\"items\");\nawait context.sync();\n\nconst highlightRange = highlightMatches.items[0];\nhighlightRange.font.highlightColor = \"Yellow\";\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the existing \"_GoBack\" bookmark (it currently sits right after\n#    \"9.13\" in the \"Chapter 9\" paragraph) so it can be re-created at its new\n#    location.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# 2) Split the \"#8.16\" run into \"#8.\" + \"16\" by inserting the \"_GoBack\"\n#    bookmark right before the \"16\".\n$findRange = $d.Content\n$findRange.Find.Execute(\"#8.16\") | Out-Null\n$splitPoint = $d.Range($findRange.End - 2, $findRange.End - 2)\n$d.Bookmarks.Add(\"_GoBack\", $splitPoint) | Out-Null\n\n# 3) Highlight \", # 8.20\" (the four runs: \",\", \" #\", \" \", \"8.20\") in yellow.\n$highlightRange = $d.Content\n$highlightRange.Find.Execute(\", # 8.20\") | Out-Null\n$highlightRange.Font.HighlightColorIndex = 7\n"}
